$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 705; existing rows 705:820 shift down to 706:821
$ws.Rows.Item(705).Insert()

# Populate the new row 705 with data
$ws.Cells.Item(705, 1).Value = 9
$ws.Cells.Item(705, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(705, 3).Value = "Metropolitana"
$ws.Cells.Item(705, 4).Value = 45218
$ws.Cells.Item(705, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(705, 5).Value = 13
$ws.Cells.Item(705, 6).Value = 100112028
$ws.Cells.Item(705, 7).Value = "Sandia"
$ws.Cells.Item(705, 8).Value = "Sin especificar"
$ws.Cells.Item(705, 9).Value = "Primera"
$ws.Cells.Item(705, 10).Value = 820
$ws.Cells.Item(705, 11).Value = 550
$ws.Cells.Item(705, 12).Value = 600
$ws.Cells.Item(705, 13).Value = 584
$ws.Cells.Item(705, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(705, 15).Value = "Perú"
$ws.Cells.Item(705, 16).Value = 584
$ws.Cells.Item(705, 17).Value = 1
$ws.Cells.Item(705, 18).Value = "Hortaliza"
